$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1) ---
# Old: A Hoje | B Operadora | C Data da Notificação | D Demanda | E Protocolo |
#      F Beneficiário | G Prazo | H Respondido | I Natureza | J Opções
# New: A Operadora | B Hoje | C Notificação | D Demanda | E Protocolo |
#      F Beneficiário | G CPF | H Descrição | I Prazo | J Respondido | K Natureza

$ws.Range("A1").Value = "Operadora"
$ws.Range("B1").Value = "Hoje"
$ws.Range("C1").Value = "Notificação"
$ws.Range("D1").Value = "Demanda"
$ws.Range("E1").Value = "Protocolo"
$ws.Range("F1").Value = "Beneficiário"
$ws.Range("G1").Value = "CPF"
$ws.Range("H1").Value = "Descrição"
$ws.Range("I1").Value = "Prazo"
$ws.Range("J1").Value = "Respondido"
$ws.Range("K1").Value = "Natureza"

# Copy header style from an existing header cell (A1) to the new K1 header
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows ---
$ws.Range("A2").Value = "417823 - PREMIUM SAÚDE S.A"
$ws.Range("B2").Value = "17-04-2023"
$ws.Range("C2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C2").Value = Get-Date -Year 2023 -Month 4 -Day 11 -Hour 8 -Minute 39 -Second 19
$ws.Range("D2").Value = 12163407
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "8600473"
$ws.Range("F2").Value = "AYLA ALVES COELHO"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "19014458606"
$ws.Range("H2").Value = "minha filha usa NEOCATE desde os 4 meses de idade, foi solicitado o pedido de intolerancia a lactose no dia 22/3/2023 porem nao tenho retorno quando ligo dizem, que pode levar ate 21 dias uteis mas em consulta ao site da ANS esse prazo seria pra procedimentos PAC, que nao e o caso de um exame de sangue para detectar intolerância a lactose, preciso de um retorno visto que o prazo de 10 dias uteis finalizou em 06/04/2023."
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = "NO"
$ws.Range("K2").Value = "Assistencial"

$ws.Range("A3").Value = "417823 - PREMIUM SAÚDE S.A"
$ws.Range("B3").Value = "17-04-2023"
$ws.Range("C3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C3").Value = Get-Date -Year 2023 -Month 4 -Day 11 -Hour 11 -Minute 9 -Second 14
$ws.Range("D3").Value = 12163869
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "8601052"
$ws.Range("F3").Value = "PRISCILA APARECIDA SANTOS FRANCISCO"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "10136083641"
$ws.Range("H3").Value = "Solicitei com 10 dias de antecedência a marcação do exame de ultrassonografia endovaginal,na clínica Santa Helena ltda através do plano,ao chegar no estabelecimento no dia do exame a operadora do plano negou meu procedimento alegando que a clinica não estava mais cadastrada,porém não recebi nenhum contanto prévio mesmo estando com o agendamento feito a dias."
$ws.Range("I3").Value = 7
$ws.Range("J3").Value = "NO"
$ws.Range("K3").Value = "Assistencial"
